$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CMS")

# Row 4 - Inventory
$ws.Range("C4").Value = 576000000.0
$ws.Range("D4").Value = 639000000.0
$ws.Range("E4").Value = 506000000.0
$ws.Range("F4").Value = 434000000.0

# Row 14 - Accounts Payable
$ws.Range("C14").Value = 678000000.0
$ws.Range("D14").Value = 667000000.0
$ws.Range("E14").Value = 635000000.0
$ws.Range("F14").Value = 497000000.0

# Row 24 - Long Term Tax Liability (Deferred)
$ws.Range("C24").Value = 1978000000.0
$ws.Range("D24").Value = 1937000000.0
$ws.Range("E24").Value = 1881000000.0
$ws.Range("F24").Value = 1841000000.0

# Row 36 - Net Debt
$ws.Range("B36").Value = 16100000000.0
$ws.Range("G36").Value = 13000000000.0

# Row 37 - Total Debt
$ws.Range("B37").Value = 16624000000.0
$ws.Range("G37").Value = 13157000000.0
